$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("G1").Value = "Helper strain"
$ws.Range("H1").Value = "Form colonies on carbon-free agar plate"

# --- Data values for new columns G (Helper strain) and H (colony growth) ---
$gValues = @{
    2  = "N/A";  3  = "N/A";  4  = "Possible"; 5  = "N/A";  6  = "Possible"
    7  = "N/A";  8  = "N/A";  9  = "N/A";       10 = "N/A";  11 = "N/A"
    12 = "N/A";  13 = "N/A";  14 = "N/A";       15 = "N/A";  16 = "N/A"
    17 = "N/A";  18 = "N/A";  19 = "N/A";       20 = "N/A";  21 = "N/A"
    22 = "N/A";  23 = "N/A";  24 = "N/A";       25 = "N/A";  26 = "N/A"
    27 = "N/A";  28 = "N/A";  29 = "N/A"
}

$hValues = @{
    2  = "Many"; 3  = "Many"; 4  = "Few";      5  = "Many"; 6  = "Few"
    7  = "Many"; 8  = "Many"; 9  = "Many";      10 = "Many"; 11 = "Very few"
    12 = "Very few"; 13 = "Yes"; 14 = "Yes";    15 = "Few";  16 = "Yes"
    17 = "Yes";  18 = "Yes";  19 = "Few";       20 = "Yes";  21 = "Yes"
    22 = "Yes";  23 = "Yes";  24 = "Yes";       25 = "Few";  26 = "Yes"
    27 = "Few";  28 = "Yes";  29 = "Very few"
}

for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 7).Value = $gValues[$r]
    $ws.Cells.Item($r, 8).Value = $hValues[$r]
}

# --- Copy formatting (style index used by column A/B/C, e.g. A2) onto the
#     new G/H cells for every row except 4 and 6, which stay unstyled,
#     matching the original authored workbook. Done in contiguous blocks
#     since multi-area paste only affects the first area in this engine. ---
$styleBlocks = @(
    @(2, 3),
    @(5, 5),
    @(7, 29)
)

foreach ($block in $styleBlocks) {
    $first = $block[0]
    $last = $block[1]
    $ws.Range("A2").Copy()
    $ws.Range("G$first`:H$last").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Restore the active selection recorded by Excel after editing ---
$ws.Range("E24").Select()
